$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (post-update TPM-derived statistics) for the Fn1-Tshr LR-pair sheet.
# Keyed by cell address -> new numeric value, taken from the updated natmi output.
$updates = @{
    "G2" = 29.223446
    "H2" = 87.670338
    "I2" = 0.0169041244192178
    "J2" = 0.0169041244192178
    "K2" = 3
    "L2" = 1
    "M2" = 0.6327629999999999
    "N2" = 1.898289
    "O2" = 0.1382544270550543
    "P2" = 0.1382544270550544
    "Q2" = 18.491515361298
    "R2" = 166.423638251682
    "S2" = 0.00233707003644631
    "T2" = 0.002337070036446311
    "G3" = 29.223446
    "H3" = 87.670338
    "I3" = 0.0169041244192178
    "J3" = 0.0169041244192178
    "O3" = 0.4765301499162115
    "P3" = 0.4765301499162115
    "Q3" = 63.73585840971533
    "R3" = 573.622725687438
    "S3" = 0.00805532494369215
    "T3" = 0.008055324943692151
    "G4" = 29.223446
    "H4" = 87.670338
    "I4" = 0.0169041244192178
    "J4" = 0.0169041244192178
    "M4" = 1.444396333333334
    "N4" = 4.333189000000001
    "O4" = 0.3155908096798033
    "P4" = 0.3155908096798033
    "Q4" = 42.21023824976467
    "R4" = 379.8921442478821
    "S4" = 0.005334786312389081
    "T4" = 0.005334786312389081
    "G5" = 29.223446
    "H5" = 87.670338
    "I5" = 0.0169041244192178
    "J5" = 0.0169041244192178
    "K5" = 3
    "L5" = 1
    "M5" = 0.3186579999999999
    "N5" = 0.9559739999999999
    "O5" = 0.06962461334893082
    "P5" = 0.06962461334893082
    "Q5" = 9.312284855467999
    "R5" = 83.810563699212
    "S5" = 0.001176943126690259
    "T5" = 0.001176943126690259
    "I6" = 0.9471112884046843
    "J6" = 0.9471112884046842
    "M6" = 0.6327629999999999
    "N6" = 1.898289
    "O6" = 0.1382544270550543
    "P6" = 0.1382544270550544
    "Q6" = 1036.05028595763
    "R6" = 9324.452573618668
    "S6" = 0.130942328535764
    "T6" = 0.130942328535764
    "I7" = 0.9471112884046843
    "J7" = 0.9471112884046842
    "O7" = 0.4765301499162115
    "P7" = 0.4765301499162115
    "S7" = 0.4513270842508204
    "T7" = 0.4513270842508204
    "I8" = 0.9471112884046843
    "J8" = 0.9471112884046842
    "M8" = 1.444396333333334
    "N8" = 4.333189000000001
    "O8" = 0.3155908096798033
    "P8" = 0.3155908096798033
    "Q8" = 2364.972721518408
    "R8" = 21284.75449366567
    "S8" = 0.298899618364516
    "T8" = 0.298899618364516
    "I9" = 0.9471112884046843
    "J9" = 0.9471112884046842
    "K9" = 3
    "L9" = 1
    "M9" = 0.3186579999999999
    "N9" = 0.9559739999999999
    "O9" = 0.06962461334893082
    "P9" = 0.06962461334893082
    "Q9" = 521.7525550999133
    "R9" = 4695.772995899219
    "S9" = 0.06594225725358385
    "T9" = 0.06594225725358384
    "G10" = 37.39212666666667
    "H10" = 112.17638
    "I10" = 0.02162924801792661
    "J10" = 0.0216292480179266
    "M10" = 0.6327629999999999
    "N10" = 1.898289
    "O10" = 0.1382544270550543
    "P10" = 0.1382544270550544
    "Q10" = 23.66035424598
    "R10" = 212.94318821382
    "S10" = 0.002990339292350113
    "T10" = 0.002990339292350113
    "G11" = 37.39212666666667
    "H11" = 112.17638
    "I11" = 0.02162924801792661
    "J11" = 0.0216292480179266
    "O11" = 0.4765301499162115
    "P11" = 0.4765301499162115
    "Q11" = 81.55161752193112
    "R11" = 733.9645576973801
    "S11" = 0.01030698880055749
    "T11" = 0.01030698880055749
    "G12" = 37.39212666666667
    "H12" = 112.17638
    "I12" = 0.02162924801792661
    "J12" = 0.0216292480179266
    "M12" = 1.444396333333334
    "N12" = 4.333189000000001
    "O12" = 0.3155908096798033
    "P12" = 0.3155908096798033
    "Q12" = 54.0090506528689
    "R12" = 486.0814558758202
    "S12" = 0.006825991894742739
    "T12" = 0.006825991894742738
    "G13" = 37.39212666666667
    "H13" = 112.17638
    "I13" = 0.02162924801792661
    "J13" = 0.0216292480179266
    "K13" = 3
    "L13" = 1
    "M13" = 0.3186579999999999
    "N13" = 0.9559739999999999
    "O13" = 0.06962461334893082
    "P13" = 0.06962461334893082
    "Q13" = 11.91530029934667
    "R13" = 107.23770269412
    "S13" = 0.001505928030276268
    "T13" = 0.001505928030276268
    "G14" = 24.817167
    "H14" = 74.45150100000001
    "I14" = 0.01435533915817136
    "J14" = 0.01435533915817136
    "M14" = 0.6327629999999999
    "N14" = 1.898289
    "O14" = 0.1382544270550543
    "P14" = 0.1382544270550544
    "Q14" = 15.703385042421
    "R14" = 141.330465381789
    "S14" = 0.001984689190493967
    "T14" = 0.001984689190493968
    "G15" = 24.817167
    "H15" = 74.45150100000001
    "I15" = 0.01435533915817136
    "J15" = 0.01435533915817136
    "O15" = 0.4765301499162115
    "P15" = 0.4765301499162115
    "Q15" = 54.125835879939
    "R15" = 487.1325229194511
    "S15" = 0.00684075192114146
    "T15" = 0.00684075192114146
    "G16" = 24.817167
    "H16" = 74.45150100000001
    "I16" = 0.01435533915817136
    "J16" = 0.01435533915817136
    "M16" = 1.444396333333334
    "N16" = 4.333189000000001
    "O16" = 0.3155908096798033
    "P16" = 0.3155908096798033
    "Q16" = 35.84582501852101
    "R16" = 322.6124251666891
    "S16" = 0.004530413108155486
    "T16" = 0.004530413108155485
    "G17" = 24.817167
    "H17" = 74.45150100000001
    "I17" = 0.01435533915817136
    "J17" = 0.01435533915817136
    "K17" = 3
    "L17" = 1
    "M17" = 0.3186579999999999
    "N17" = 0.9559739999999999
    "O17" = 0.06962461334893082
    "P17" = 0.06962461334893082
    "Q17" = 7.908188801885999
    "R17" = 71.173699216974
    "S17" = 0.000999484938380447
    "T17" = 0.0009994849383804468
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

